# Fix "4 T's" -> "5 T's" in the "Reversible Causes (7 H's and 4 T's)" table
# header, splitting the run that holds " (7 H's and 4 T's)" into three runs
# (" (7 H's and ", "5", " T's)") the same way the authored edit did, while
# leaving the preceding "Reversible Causes" run untouched.

$d = $word.ActiveDocument

$rsquo = [char]0x2019

# Locate the run of text containing the digit we need to change.
$ctx = $d.Content
$found = $ctx.Find.Execute("7 H" + $rsquo + "s and 4 T" + $rsquo + "s)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the '7 H's and 4 T's' text to update."
}

$digitStart = $ctx.Start + $ctx.Text.IndexOf("4")
$digitEnd = $digitStart + 1
$digitRange = $d.Range($digitStart, $digitEnd)

# Toggling a character property (without touching the text) splits the
# " (7 H's and 4 T's)" run into three runs around the "4" without merging
# the preceding "Reversible Causes" run into it.
$digitRange.Bold = 1

# While the "4" run is still distinctly formatted (Bold) from its
# neighbours, swap its text for "5" - this keeps the three-way split
# intact instead of re-coalescing the surrounding runs.
$digitRange.Characters(1).Text = "5"

# Re-select the (still single-character) range and drop the temporary
# Bold flag; this is a property change, not a text change, so it does not
# trigger another run merge.
$newDigitRange = $d.Range($digitStart, $digitStart + 1)
$newDigitRange.Bold = 0
